# Auto-generated PowerShell (Excel COM-interop) script
# Applies cached numeric value updates to the Goblin Profits workbook sheets
# (columns H-N: currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 96.40000000000001
$ws.Range("I11").Value = 96.40000000000001
$ws.Range("K11").Value = 96.40000000000001
$ws.Range("M11").Value = 43.59999999999999
$ws.Range("H17").Value = 1982.56
$ws.Range("J17").Value = 2180.05
$ws.Range("L17").Value = 6540.150000000001
$ws.Range("N17").Value = -6876.150000000001
$ws.Range("H19").Value = 466.33334
$ws.Range("I19").Value = 599.5
$ws.Range("K19").Value = 599.5
$ws.Range("M19").Value = -424.5
$ws.Range("H40").Value = 3499.3333
$ws.Range("I40").Value = 2600
$ws.Range("J40").Value = 3949
$ws.Range("K40").Value = 2600
$ws.Range("L40").Value = 3949
$ws.Range("M40").Value = -2425
$ws.Range("N40").Value = -4299
$ws.Range("H92").Value = 642.4211
$ws.Range("I92").Value = 549.88
$ws.Range("K92").Value = 549.88
$ws.Range("M92").Value = 698.12
$ws.Range("H100").Value = 2968.375
$ws.Range("J100").Value = 3407.2
$ws.Range("L100").Value = 3407.2
$ws.Range("N100").Value = -4489.2
$ws.Range("H112").Value = 4125
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4125
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 12375
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -14591
$ws.Range("H113").Value = 3798.8462
$ws.Range("I113").Value = 3439.5
$ws.Range("K113").Value = 3439.5
$ws.Range("M113").Value = -185.5
$ws.Range("H132").Value = 2582.3076
$ws.Range("I132").Value = 1779.409
$ws.Range("K132").Value = 5338.227000000001
$ws.Range("M132").Value = -2808.227000000001
$ws.Range("H137").Value = 1722.5
$ws.Range("I137").Value = 1503.0646
$ws.Range("J137").Value = 2122.647
$ws.Range("K137").Value = 4509.1938
$ws.Range("L137").Value = 6367.941
$ws.Range("M137").Value = -1959.1938
$ws.Range("N137").Value = -11467.941
$ws.Range("H138").Value = 7775.398
$ws.Range("I138").Value = 6619.9473
$ws.Range("J138").Value = 8093.5654
$ws.Range("K138").Value = 19859.8419
$ws.Range("L138").Value = 24280.6962
$ws.Range("M138").Value = -14719.8419
$ws.Range("N138").Value = -34560.69620000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1411.5834
$ws.Range("I2").Value = 1493.8889
$ws.Range("K2").Value = 1493.8889
$ws.Range("M2").Value = -1380.8889
$ws.Range("H6").Value = 2499
$ws.Range("I6").Value = 2499
$ws.Range("K6").Value = 2499
$ws.Range("M6").Value = -2326
$ws.Range("H32").Value = 4248.35
$ws.Range("I32").Value = 3357.75
$ws.Range("J32").Value = 12263.75
$ws.Range("K32").Value = 3357.75
$ws.Range("L32").Value = 12263.75
$ws.Range("M32").Value = -3070.75
$ws.Range("N32").Value = -12837.75
$ws.Range("H61").Value = 1433.7765
$ws.Range("I61").Value = 809.56757
$ws.Range("K61").Value = 809.56757
$ws.Range("M61").Value = -597.56757
$ws.Range("H74").Value = 1824.7317
$ws.Range("J74").Value = 4881.6
$ws.Range("L74").Value = 4881.6
$ws.Range("N74").Value = -6629.6
$ws.Range("H77").Value = 1824.7317
$ws.Range("J77").Value = 4881.6
$ws.Range("L77").Value = 24408
$ws.Range("N77").Value = -33144
$ws.Range("H80").Value = 74999.5
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 74999.5
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H102").Value = 6729.9473
$ws.Range("I102").Value = 4788.5557
$ws.Range("J102").Value = 8477.200000000001
$ws.Range("K102").Value = 4788.5557
$ws.Range("L102").Value = 8477.200000000001
$ws.Range("M102").Value = -3166.5557
$ws.Range("N102").Value = -11721.2
$ws.Range("H116").Value = 1411.5834
$ws.Range("I116").Value = 1493.8889
$ws.Range("K116").Value = 1493.8889
$ws.Range("M116").Value = 800.1111000000001
$ws.Range("H136").Value = 1433.7765
$ws.Range("I136").Value = 809.56757
$ws.Range("K136").Value = 2428.70271
$ws.Range("M136").Value = 121.29729

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1411.5834
$ws.Range("I3").Value = 1493.8889
$ws.Range("K3").Value = 1493.8889
$ws.Range("M3").Value = -1379.8889
$ws.Range("H86").Value = 20840136
$ws.Range("I86").Value = 10715.25
$ws.Range("J86").Value = 41669556
$ws.Range("K86").Value = 10715.25
$ws.Range("L86").Value = 41669556
$ws.Range("M86").Value = -9592.25
$ws.Range("N86").Value = -41671802
$ws.Range("H89").Value = 20840136
$ws.Range("I89").Value = 10715.25
$ws.Range("J89").Value = 41669556
$ws.Range("K89").Value = 53576.25
$ws.Range("L89").Value = 208347780
$ws.Range("M89").Value = -47960.25
$ws.Range("N89").Value = -208359012
$ws.Range("H99").Value = 3179
$ws.Range("I99").Value = 2124.8333
$ws.Range("J99").Value = 4233.1665
$ws.Range("K99").Value = 2124.8333
$ws.Range("L99").Value = 4233.1665
$ws.Range("M99").Value = -626.8332999999998
$ws.Range("N99").Value = -7229.1665
$ws.Range("H105").Value = 33720
$ws.Range("I105").Value = 33720
$ws.Range("K105").Value = 33720
$ws.Range("M105").Value = -31973

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3533.2258
$ws.Range("I7").Value = 325.76923
$ws.Range("K7").Value = 325.76923
$ws.Range("M7").Value = -212.76923
$ws.Range("H16").Value = 1199.2222
$ws.Range("I16").Value = 974.125
$ws.Range("K16").Value = 974.125
$ws.Range("M16").Value = -687.125
$ws.Range("H31").Value = 3388.6875
$ws.Range("I31").Value = 1657
$ws.Range("J31").Value = 5615.143
$ws.Range("K31").Value = 1657
$ws.Range("L31").Value = 5615.143
$ws.Range("M31").Value = -1362
$ws.Range("N31").Value = -6205.143
$ws.Range("H34").Value = 3388.6875
$ws.Range("I34").Value = 1657
$ws.Range("J34").Value = 5615.143
$ws.Range("K34").Value = 1657
$ws.Range("L34").Value = 5615.143
$ws.Range("M34").Value = -1455
$ws.Range("N34").Value = -6019.143
$ws.Range("H105").Value = 4296
$ws.Range("I105").Value = 5506.7144
$ws.Range("J105").Value = 3085.2856
$ws.Range("K105").Value = 5506.7144
$ws.Range("L105").Value = 3085.2856
$ws.Range("M105").Value = -3759.7144
$ws.Range("N105").Value = -6579.2856
$ws.Range("H107").Value = 662.38464
$ws.Range("I107").Value = 524.1579
$ws.Range("J107").Value = 1037.5714
$ws.Range("K107").Value = 524.1579
$ws.Range("L107").Value = 1037.5714
$ws.Range("M107").Value = 1395.8421
$ws.Range("N107").Value = -4877.5714
$ws.Range("H113").Value = 1199.2222
$ws.Range("I113").Value = 974.125
$ws.Range("K113").Value = 974.125
$ws.Range("M113").Value = 1195.875
$ws.Range("H122").Value = 1238.5625
$ws.Range("J122").Value = 1362.2222
$ws.Range("L122").Value = 4086.6666
$ws.Range("N122").Value = -8986.6666
$ws.Range("H132").Value = 1214.5862
$ws.Range("I132").Value = 1127.3462
$ws.Range("J132").Value = 1970.6666
$ws.Range("K132").Value = 3382.0386
$ws.Range("L132").Value = 5911.9998
$ws.Range("M132").Value = -852.0385999999999
$ws.Range("N132").Value = -10971.9998
$ws.Range("H141").Value = 125712.43
$ws.Range("J141").Value = 125712.43
$ws.Range("L141").Value = 125712.43
$ws.Range("N141").Value = -136072.43

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 245348.89
$ws.Range("I11").Value = 308592.7
$ws.Range("J11").Value = 2914.3333
$ws.Range("K11").Value = 925778.1000000001
$ws.Range("L11").Value = 8742.999899999999
$ws.Range("M11").Value = -925638.1000000001
$ws.Range("N11").Value = -9022.999899999999
$ws.Range("H57").Value = 6136.364
$ws.Range("I57").Value = 2250
$ws.Range("K57").Value = 6750
$ws.Range("M57").Value = -6191
$ws.Range("H97").Value = 143.55556
$ws.Range("I97").Value = 97.5
$ws.Range("J97").Value = 156.71428
$ws.Range("K97").Value = 292.5
$ws.Range("L97").Value = 470.14284
$ws.Range("M97").Value = 203.5
$ws.Range("N97").Value = -1462.14284
$ws.Range("H103").Value = 1072.7778
$ws.Range("I103").Value = 253.75
$ws.Range("K103").Value = 761.25
$ws.Range("M103").Value = 117.75
$ws.Range("H112").Value = 4283.8184
$ws.Range("I112").Value = 2999
$ws.Range("K112").Value = 8997
$ws.Range("M112").Value = -7889
$ws.Range("H118").Value = 3543.9092
$ws.Range("J118").Value = 5666.6665
$ws.Range("L118").Value = 16999.9995
$ws.Range("N118").Value = -19485.9995
$ws.Range("H131").Value = 2265.0527
$ws.Range("I131").Value = 786.7619
$ws.Range("K131").Value = 2360.2857
$ws.Range("M131").Value = 2679.7143

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 984.2558
$ws.Range("J97").Value = 698.55554
$ws.Range("L97").Value = 698.55554
$ws.Range("N97").Value = -1690.55554
$ws.Range("H113").Value = 52641052
$ws.Range("I113").Value = 71438030
$ws.Range("K113").Value = 71438030
$ws.Range("M113").Value = -71435860
$ws.Range("H132").Value = 1578.4
$ws.Range("I132").Value = 1723
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5169
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2639
$ws.Range("N132").Value = -8060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4468.4443
$ws.Range("I16").Value = 591.0625
$ws.Range("K16").Value = 591.0625
$ws.Range("M16").Value = -421.0625
$ws.Range("H22").Value = 3411.7
$ws.Range("I22").Value = 2718
$ws.Range("J22").Value = 3874.1667
$ws.Range("K22").Value = 2718
$ws.Range("L22").Value = 3874.1667
$ws.Range("M22").Value = -2423
$ws.Range("N22").Value = -4464.1667
$ws.Range("H27").Value = 3411.7
$ws.Range("I27").Value = 2718
$ws.Range("J27").Value = 3874.1667
$ws.Range("K27").Value = 2718
$ws.Range("L27").Value = 3874.1667
$ws.Range("M27").Value = -2611
$ws.Range("N27").Value = -4088.1667
$ws.Range("H55").Value = 1346.15
$ws.Range("I55").Value = 164
$ws.Range("K55").Value = 164
$ws.Range("M55").Value = 9
$ws.Range("H61").Value = 5462.533
$ws.Range("I61").Value = 5448.5264
$ws.Range("J61").Value = 5486.727
$ws.Range("K61").Value = 5448.5264
$ws.Range("L61").Value = 5486.727
$ws.Range("M61").Value = -5246.5264
$ws.Range("N61").Value = -5890.727
$ws.Range("H100").Value = 3332.3
$ws.Range("I100").Value = 4181
$ws.Range("J100").Value = 2766.5
$ws.Range("K100").Value = 4181
$ws.Range("L100").Value = 2766.5
$ws.Range("M100").Value = -3640
$ws.Range("N100").Value = -3848.5
$ws.Range("H113").Value = 5462.533
$ws.Range("I113").Value = 5448.5264
$ws.Range("J113").Value = 5486.727
$ws.Range("K113").Value = 5448.5264
$ws.Range("L113").Value = 5486.727
$ws.Range("M113").Value = -3278.5264
$ws.Range("N113").Value = -9826.726999999999
$ws.Range("H132").Value = 3561.2222
$ws.Range("I132").Value = 3127.7585
$ws.Range("K132").Value = 9383.2755
$ws.Range("M132").Value = -6853.2755
$ws.Range("H136").Value = 1644.8507
$ws.Range("I136").Value = 1468
$ws.Range("K136").Value = 4404
$ws.Range("M136").Value = -1854

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 38236
$ws.Range("J95").Value = 38236
$ws.Range("L95").Value = 38236
$ws.Range("N95").Value = -43728
$ws.Range("H136").Value = 1204.5172
$ws.Range("I136").Value = 791.14
$ws.Range("K136").Value = 2373.42
$ws.Range("M136").Value = 176.5799999999999

